$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column H/I data mirroring columns A/B with some page numbers updated.
$ws.Range("H1").Value = "低龄："
$ws.Range("I1").Value = "高龄："

$ws.Range("H2").Value = "page1"
$ws.Range("I2").Value = "page1"

$ws.Range("H3").Value = "page2"
$ws.Range("I3").Value = "page2"

$ws.Range("H4").Value = "page6"
$ws.Range("I4").Value = "page6"

$ws.Range("H5").Value = "page50"
$ws.Range("I5").Value = "page50"

$ws.Range("H6").Value = "page4"
$ws.Range("I6").Value = "page4"

$ws.Range("H7").Value = "page5"
$ws.Range("I7").Value = "page5"

$ws.Range("H8").Value = "page51"
$ws.Range("I8").Value = "page52"

$ws.Range("H9").Value = "page10"
$ws.Range("I9").Value = "page32"

$ws.Range("H10").Value = "page11"
$ws.Range("I10").Value = "page33"

$ws.Range("H11").Value = "page13"
$ws.Range("I11").Value = "page55"

$ws.Range("H12").Value = "page14"

$ws.Range("H13").Value = "page15"
$ws.Range("I13").Value = "page15"

$ws.Range("H14").Value = "page16"
$ws.Range("I14").Value = "page35"

$ws.Range("H15").Value = "page17"

$ws.Range("H16").Value = "page18"
$ws.Range("I16").Value = "page18"

$ws.Range("H17").Value = "page22"
$ws.Range("I17").Value = "page22"

$ws.Range("H18").Value = "page23"
$ws.Range("I18").Value = "page23"

$ws.Range("H19").Value = "page53"
$ws.Range("I19").Value = "page54"

$ws.Range("H20").Value = "page26"
$ws.Range("I20").Value = "page26"

# Update the selection to reflect the newly filled column I.
$ws.Range("I2:I20").Select()
